$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author selected A19:I28 and cleared the contents of the now-unneeded
# footnote / legend rows (rows 19-26), leaving the row/column structure and
# cell formatting in place but removing the text.
$ws.Range("A19:I28").ClearContents()

$ws.Range("I28").Select()
